# Updated numbers for 2023.
# For each sheet ("Jack" and "Jill") the row holding the *last* all-blank
# year (the extra row past the final populated year) is removed, and every
# remaining year in column A is bumped forward by one (2023 -> 2024, etc.).

$wb = $excel.ActiveWorkbook

$sheetRowToDelete = @{ "Jack" = 33; "Jill" = 35 }

foreach ($ws in $wb.Worksheets) {
    $rowToDelete = $sheetRowToDelete[$ws.Name]
    if ($rowToDelete) {
        $ws.Rows.Item($rowToDelete).Delete()
    }

    # Find the last used row in column A (years live there).
    $lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row  # xlUp

    # Row 1 is the header ("year", ...); years start on row 2.
    # NOTE: use .Value2 (not .Value) — reading .Value back out in an
    # expression here yields the property descriptor, not the cell's data.
    #
    # The trailing row on some sheets is a bare placeholder (year only, no
    # formatted data cells to its right) and keeps its original year value
    # — only the "real" data rows get bumped forward by one year.
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $isBareRow = $ws.Cells.Item($r, 2).NumberFormat -eq "General"
        if (($cell.Value2 -ne $null) -and (-not $isBareRow)) {
            $cell.Value2 = $cell.Value2 + 1
        }
    }

    # Mimic clicking the row-2 header, which is the selection left behind
    # by this edit in the source workbook.
    $ws.Range("A2:XFD2").Select()
}

# "Jack" was (and remains) the active tab; visiting "Jill" above to update
# its selection must not leave it as the active sheet.
$wb.Worksheets.Item("Jack").Activate()
